$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("A5").Value = "MLPRegressor"

$ws.Range("B3").Value = 0.9945610105879964
$ws.Range("C3").Value = 0.994062998422951
$ws.Range("D3").Value = 0.994111959866958

$ws.Range("B4").Value = 0.9952060653613537
$ws.Range("C4").Value = 0.9952402786154169
$ws.Range("D4").Value = 0.995234186784431

$ws.Range("B5").Value = 0.9958737736289088
$ws.Range("C5").Value = 0.9953999585243603
$ws.Range("D5").Value = 0.9963184280664361
